# Auto-generated Excel COM-interop edit script
$wb = $excel.ActiveWorkbook

# --- Sheet: investigaciones ---
$wsInv = $wb.Worksheets.Item("investigaciones")
$wsInv.Cells.Item(34, 8).Value = '5 Finalizada'

# --- Sheet: productos ---
$wsProd = $wb.Worksheets.Item("productos")
$wsProd.Cells.Item(550, 2).Value = 'Visualización'
$wsProd.Cells.Item(550, 3).Value = 'Tablero de resultados cuantitativos'
$wsProd.Cells.Item(550, 4).Value = 'No'
$wsProd.Cells.Item(550, 5).Value = 'https://lookerstudio.google.com/u/0/reporting/4a928fd5-3d5a-4de3-a655-943ad2da4e6e/page/p_d6cfr52jwd'
$wsProd.Cells.Item(574, 2).Value = 'Visualización'
$wsProd.Cells.Item(574, 3).Value = 'Tablero de resultados cuantitativos'
$wsProd.Cells.Item(574, 4).Value = 'Sí'
$wsProd.Cells.Item(574, 5).Value = 'https://lookerstudio.google.com/reporting/f0ca73ce-6555-49ec-9219-1ee6b4e60428'
$wsProd.Cells.Item(578, 4).Value = 'Sí'
$wsProd.Cells.Item(579, 5).Value = 'https://drive.google.com/file/d/1zSXYH6cCKPEoIcuur0PuPncZ0IaOLo8E/view?usp=drive_link'
$wsProd.Cells.Item(581, 2).Value = 'Presentación'
$wsProd.Cells.Item(581, 3).Value = 'Presentación resultados Festival Monumentum 2025'
$wsProd.Cells.Item(581, 4).Value = 'Sí'
$wsProd.Cells.Item(581, 5).Value = 'https://drive.google.com/file/d/1xJO1lFA1uh52DD5A-mAYv6s-glGlN2jZ/view?usp=sharing'
$wsProd.Cells.Item(604, 3).Value = 'Anexo 1. Análisis de las preguntas espaciales'
$wsProd.Cells.Item(605, 3).Value = 'Productos finales Festival Joropo al Parque'
$wsProd.Cells.Item(605, 5).Value = 'https://drive.google.com/drive/folders/1Lcu0iYLFoiXHMX8f8ivPgJPF5XDn-IYx?usp=drive_link'
$wsProd.Cells.Item(606, 3).Value = 'Productos finales Festival Rock al Parque'
$wsProd.Cells.Item(606, 5).Value = 'https://drive.google.com/drive/folders/1rccxA6SFUSsKbjgRKqWnvEbok0SpWS9t?usp=drive_link'
$wsProd.Cells.Item(607, 3).Value = 'Productos finales Festival Vallenato al Parque'
$wsProd.Cells.Item(607, 5).Value = 'https://drive.google.com/drive/folders/1BchZfxR7zuYYbkgIH8qE8fdgQ1y1EYZk?usp=drive_link'
$wsProd.Cells.Item(608, 3).Value = 'Productos finales Festival Colombia al Parque'
$wsProd.Cells.Item(608, 5).Value = 'https://drive.google.com/drive/folders/1HBmOWH-hQDqwj0mDiTDqlVdSUVcQkhYt?usp=drive_link'
$wsProd.Cells.Item(609, 3).Value = 'Productos finales Festival Jazz al Parque'
$wsProd.Cells.Item(609, 5).Value = 'https://drive.google.com/drive/folders/1SSz2sXnITpz_lT4g5-n2AR3M8Ed3LsnM?usp=drive_link'
$wsProd.Cells.Item(616, 3).Value = 'Productos finales Festival Popular al Parque'
$wsProd.Cells.Item(616, 5).Value = 'https://drive.google.com/drive/folders/14CmZAQrG8eHOh7xDTf8ypS7TweGA5Df_?usp=drive_link'
$wsProd.Cells.Item(617, 3).Value = 'Productos finales Festival Hip Hop al Parque'
$wsProd.Cells.Item(617, 5).Value = 'https://drive.google.com/drive/folders/19PXuTzH3C_Nfl5Gj8LO1_XHUgBfjb4b4?usp=drive_link'
$wsProd.Cells.Item(618, 3).Value = 'Productos finales Festival Salsa al Parque'
$wsProd.Cells.Item(618, 5).Value = 'https://drive.google.com/drive/folders/1DwWfpUiZRNgalk9jBxS-hFLcAWPqdcap?usp=drive_link'
$wsProd.Cells.Item(621, 2).Value = 'Presentación'
$wsProd.Cells.Item(621, 3).Value = 'Presentación Festival Patrimonios en Ruana 2025'
$wsProd.Cells.Item(621, 4).Value = 'Sí'
$wsProd.Cells.Item(621, 5).Value = 'https://drive.google.com/file/d/1CaHBM28hdOyauLNWrACrgLpvdPbUDqch/view?usp=sharing'
$wsProd.Cells.Item(622, 1).Value = 110
$wsProd.Cells.Item(622, 2).Value = 'Presentación'
$wsProd.Cells.Item(622, 3).Value = 'Presentación Noche de Museos 2025'
$wsProd.Cells.Item(622, 4).Value = 'Sí'
$wsProd.Cells.Item(622, 5).Value = 'https://drive.google.com/file/d/1CaHBM28hdOyauLNWrACrgLpvdPbUDqch/view?usp=sharing'
$wsProd.Cells.Item(670, 5).Value = 'https://drive.google.com/file/d/1mTOdJ3SeDN8tW2UEmwBNQ1ineXgoHsou/view?usp=sharing'
$wsProd.Cells.Item(704, 3).Value = 'Informe final'
$wsProd.Cells.Item(705, 3).Value = 'Anexo de definición'
$wsProd.Cells.Item(706, 3).Value = 'Formulario sondeo'

# --- Sheet: hallazgos ---
$wsHal = $wb.Worksheets.Item("hallazgos")
$wsHal.Rows.Item(293).Resize(32).Insert()

$wsHal.Cells.Item(293, 1).Value = 104
$wsHal.Cells.Item(293, 2).Value = 1
$wsHal.Cells.Item(293, 3).Value = 'Experiencia y participación cultural'
$wsHal.Cells.Item(293, 4).Value = 'El Festival Monumentum 2025 registra una experiencia altamente valorada por sus asistentes y un desempeño sólido en calidad artística. Predominan asistentes primerizos, lo que sugiere capacidad de ampliación de públicos; al mismo tiempo, se observa una alta disposición a volver y a recomendar el evento. La participación no se limita a la asistencia: aparecen prácticas asociadas de circulación y apropiación (seguimiento de contenidos, espacios formativos y socialización de la experiencia), lo que refuerza el festival como un nodo que activa participación cultural más sostenida en Bogotá D.C. '
$wsHal.Cells.Item(294, 1).Value = 104
$wsHal.Cells.Item(294, 2).Value = 2
$wsHal.Cells.Item(294, 3).Value = 'Convivencia, identidad y orgullo por Bogotá D.C.'
$wsHal.Cells.Item(294, 4).Value = 'Los resultados muestran que el evento se percibe como un espacio de encuentro, respeto e inclusión en el espacio público, con valoraciones favorables sobre convivencia y comportamiento durante la jornada. En el plano simbólico, se reporta una conexión fuerte con la cultura de la ciudad y con el arte como forma de transformación, junto con un componente extendido de orgullo por Bogotá D.C. y por su oferta cultural. En conjunto, el festival aparece como una experiencia urbana que fortalece sentidos de pertenencia y reconocimiento cultural. '
$wsHal.Cells.Item(295, 1).Value = 104
$wsHal.Cells.Item(295, 2).Value = 3
$wsHal.Cells.Item(295, 3).Value = 'Condiciones de acceso, comunicación y sostenibilidad'
$wsHal.Cells.Item(295, 4).Value = 'En accesibilidad, la mayoría reporta facilidad para gestionar transporte de llegada, ingreso, información y aspectos logísticos generales, aunque persisten retos puntuales en seguridad y planeación del transporte de salida para una fracción del público. En mediación y comunicación, predomina el canal digital (especialmente redes), con percepciones mayoritarias de información clara, pero con margen para fortalecer la promoción. En sostenibilidad, se destaca un patrón de movilidad asociado principalmente a transporte público masivo, mientras que la percepción sobre separación de residuos es positiva pero menos contundente, lo que sugiere oportunidades de mejora en visibilidad y apropiación de acciones ambientales. '
$wsHal.Cells.Item(296, 1).Value = 104
$wsHal.Cells.Item(296, 2).Value = 4
$wsHal.Cells.Item(296, 3).Value = 'Efectos económicos indirectos y dinámica territorial heterogénea'
$wsHal.Cells.Item(296, 4).Value = 'El análisis económico descriptivo muestra dos señales principales: el gasto declarado por asistentes en rubros asociados a su participación y el reporte de incremento de ventas en una parte de los negocios de la zona de influencia. Sin embargo, estos efectos no son homogéneos: no todos los establecimientos perciben mejoras, y la evidencia de inversión adicional y contratación de personal es limitada. En síntesis, el festival activa dinámicas económicas locales de forma parcial y desigual, aportando indicios útiles para fortalecer estrategias de articulación territorial y aprovechamiento económico en futuras ediciones en Bogotá D.C.'
$wsHal.Cells.Item(297, 1).Value = 105
$wsHal.Cells.Item(297, 2).Value = 1
$wsHal.Cells.Item(297, 3).Value = 'Libro al Viento como encuentro ocasional con la lectura'
$wsHal.Cells.Item(297, 4).Value = 'Los resultados muestran que Libro al Viento funciona principalmente como una oportunidad ocasional de encuentro con la lectura, más que como un sistema de préstamo frecuente. La mayoría de las personas que toman libros accede a pocos ejemplares y lo hace de manera puntual. Esta característica es coherente con el diseño del programa y con su presencia en el espacio público, donde el acceso es abierto, espontáneo y no mediado por trámites o registros.'
$wsHal.Cells.Item(298, 1).Value = 105
$wsHal.Cells.Item(298, 2).Value = 2
$wsHal.Cells.Item(298, 3).Value = 'Experiencia de lectura, disfrute, utilidad y apropiación
'
$wsHal.Cells.Item(298, 4).Value = 'La experiencia asociada a los libros de Libro al Viento combina el disfrute con usos prácticos. La lectura aparece vinculada principalmente a la diversión, pero también al trabajo, al aprendizaje y al enriquecimiento cultural, con variaciones según el tipo de entorno. Además, los efectos simbólicos son claros entre quienes conocen el programa: orgullo por su existencia, disfrute de la lectura y reconocimiento de que leer no es solo para especialistas. Estos impactos, aunque concentrados en un grupo reducido, son intensos y significativos.
'
$wsHal.Cells.Item(299, 1).Value = 105
$wsHal.Cells.Item(299, 2).Value = 3
$wsHal.Cells.Item(299, 3).Value = 'Presencia urbana y valor cultural del programa'
$wsHal.Cells.Item(299, 4).Value = 'Libro al Viento es percibido como un proyecto cultural con presencia en múltiples espacios de la ciudad, lo que refuerza su carácter de política cultural urbana. Su circulación por bibliotecas, parques, eventos culturales, instituciones educativas y espacios de tránsito cotidiano lo posiciona como una iniciativa accesible, cercana y reconocible, que trasciende los espacios tradicionales de lectura y se integra a la vida urbana.'
$wsHal.Cells.Item(300, 1).Value = 105
$wsHal.Cells.Item(300, 2).Value = 4
$wsHal.Cells.Item(300, 3).Value = 'Circulación extendida y apropiación simbólica de los libros'
$wsHal.Cells.Item(300, 4).Value = 'La manera en que los libros circulan muestra que su valor no se agota en el punto de entrega ni en la regla de devolución. Al conservarlos, regalarlos o compartirlos, las personas expresan una fuerte apropiación simbólica. Esta circulación extendida refuerza la presencia del programa en la vida cotidiana y lo inscribe en trayectorias personales y sociales diversas, aunque también plantea retos para la disponibilidad y sostenibilidad del sistema de dispensadores.'
$wsHal.Cells.Item(301, 1).Value = 105
$wsHal.Cells.Item(301, 2).Value = 5
$wsHal.Cells.Item(301, 3).Value = 'El libro como objeto social situado: inserción comunitaria y política cultural '
$wsHal.Cells.Item(301, 4).Value = 'Las cartografías sociales evidencian que los libros del programa Libro al Viento se integran a prácticas sociales concretas y territorialmente situadas, más allá de su simple circulación física. Los libros se insertan en dinámicas comunitarias, educativas y cotidianas, adquiriendo sentidos diversos según los contextos sociales y territoriales de cada taller. En este sentido, Libro al Viento se configura como una política cultural situada, cuyos efectos dependen de las relaciones sociales que median la apropiación del libro. Las cartografías complementan el sondeo cuantitativo al mostrar cómo los libros se anclan, circulan y se resignifican en la vida cotidiana, dando lugar a circuitos múltiples y no lineales de uso. 
'
$wsHal.Cells.Item(302, 1).Value = 106
$wsHal.Cells.Item(302, 2).Value = 1
$wsHal.Cells.Item(302, 3).Value = 'Experiencia altamente valorada y con fuerte fidelización'
$wsHal.Cells.Item(302, 4).Value = 'Los resultados muestran una valoración mayoritariamente positiva del Premio Luis Caballero 2025. La calidad de los proyectos y artistas es calificada como excelente por una proporción amplia de asistentes, y se observa una muy alta disposición a recomendar y a asistir nuevamente en futuras ediciones. En conjunto, esto reafirma al premio como una plataforma sólida de circulación del arte contemporáneo, capaz de generar satisfacción y continuidad en el público. '
$wsHal.Cells.Item(303, 1).Value = 106
$wsHal.Cells.Item(303, 2).Value = 2
$wsHal.Cells.Item(303, 3).Value = 'Ampliación de audiencias y vínculo con la oferta cultural de la ciudad'
$wsHal.Cells.Item(303, 4).Value = 'El premio convoca públicos con trayectorias culturales activas, pero también funciona como puerta de entrada: una mayoría de asistentes reporta estar participando por primera vez en el evento y, antes de la visita, una proporción importante no lo conocía o tenía un conocimiento limitado. Tras la experiencia, aumenta el reconocimiento del Premio como un espacio relevante y se fortalece la disposición a participar en otros eventos culturales de Bogotá D.C., lo que sugiere un efecto de “enganche” con la oferta cultural más amplia. '
$wsHal.Cells.Item(304, 1).Value = 106
$wsHal.Cells.Item(304, 2).Value = 3
$wsHal.Cells.Item(304, 3).Value = 'Aporte simbólico y lectura crítica del arte contemporáneo'
$wsHal.Cells.Item(304, 4).Value = 'Más allá del disfrute, los asistentes señalan que el premio contribuye a valorar el arte y el patrimonio como formas de transformación, ampliar conocimientos sobre el entorno cultural y reconocerse como parte activa de la vida cultural. En el componente modular resalta que, para una mayoría, el premio aporta al debate sobre arte contemporáneo y promueve proyectos disruptivos, pero también se reconoce que puede reforzar ciertos discursos dominantes, mostrando que el público lo comprende como un espacio de discusión real, con tensiones propias del campo artístico. '
$wsHal.Cells.Item(305, 1).Value = 106
$wsHal.Cells.Item(305, 2).Value = 4
$wsHal.Cells.Item(305, 3).Value = 'Accesibilidad y mediación con percepción positiva'
$wsHal.Cells.Item(305, 4).Value = 'En términos de logística y acceso, se reportan percepciones favorables sobre señalización, acompañamiento del personal, fluidez en entradas y/o salidas, y condiciones de accesibilidad (incluida discapacidad), mientras que transporte y seguridad aparecen como factores centrales al decidir asistir. Por otro lado, los medios de comunicación sobre el evento se distribuyen entre internet, voz a voz y descubrimiento en el lugar, lo que sugiere una oportunidad para fortalecer la divulgación previa.'
$wsHal.Cells.Item(306, 1).Value = 107
$wsHal.Cells.Item(306, 2).Value = 1
$wsHal.Cells.Item(306, 3).Value = 'El circo en Bogotá: un proyecto de vida en un sector consolidado'
$wsHal.Cells.Item(306, 4).Value = 'El sector circense en Bogotá presenta un alto grado de madurez, caracterizado por la predominancia de artistas con más de diez años de experiencia, formación empírica y/o especializada, y desempeño en múltiples disciplinas y oficios. Esta configuración da cuenta de un campo complejo, diversificado y técnicamente autónomo. De manera significativa, una parte mayoritaria de los participantes reconoce el circo como su proyecto de vida y principal ocupación, reafirmando su condición de sector profesional y no meramente recreativo.'
$wsHal.Cells.Item(307, 1).Value = 107
$wsHal.Cells.Item(307, 2).Value = 2
$wsHal.Cells.Item(307, 3).Value = 'Ecosistema polivalente y móvil, anclado en el espacio público'
$wsHal.Cells.Item(307, 4).Value = 'Aunque predomina el circo contemporáneo, el social‑comunitario y el tradicional mantienen una presencia significativa, con trayectorias mixtas e hibridaciones que se desplazan entre carpas, espacio público, salas, escuelas y espacios no convencionales. El entrenamiento y la circulación se apoyan fuertemente en parques, plazas y espacios comunitarios, lo que expresa la vocación itinerante del circo, pero también la ausencia de infraestructura estable y adecuada para ensayo y presentación.'
$wsHal.Cells.Item(308, 1).Value = 107
$wsHal.Cells.Item(308, 2).Value = 3
$wsHal.Cells.Item(308, 3).Value = 'Organización fuerte, sostenibilidad frágil'
$wsHal.Cells.Item(308, 4).Value = 'El sector muestra altos niveles de organización y articulación: cerca de dos tercios de las personas encuestadas pertenecen a redes, colectivos u organizaciones, y más del 70% ha participado en mesas, sindicatos o consejos del sector cultural. Sin embargo, la mayoría se ubica en condiciones de informalidad o semi‑formalidad laboral, con ingresos esporádicos y baja protección social, lo que tensiona la posibilidad de sostener el circo como trabajo estable pese a la fuerte capacidad asociativa'
$wsHal.Cells.Item(309, 1).Value = 107
$wsHal.Cells.Item(309, 2).Value = 4
$wsHal.Cells.Item(309, 3).Value = 'Restricciones estructurales que limitan la sostenibilidad económica y la circulación'
$wsHal.Cells.Item(309, 4).Value = 'Aunque el 59% declara vivir principalmente del circo y otro grupo importante lo ejerce como actividad complementaria o mediante proyectos, las principales dificultades para sostener la práctica son el acceso a recursos públicos, la entrada a festivales y espacios de circulación, y las restricciones para usar el espacio público. Estas barreras, sumadas a la precariedad de infraestructura y la inestabilidad de los lugares de entrenamiento y presentación, impiden traducir la experiencia acumulada en condiciones económicas más estables.'
$wsHal.Cells.Item(310, 1).Value = 107
$wsHal.Cells.Item(310, 2).Value = 5
$wsHal.Cells.Item(310, 3).Value = 'Un sector formado que exige condiciones para profesionalizarse'
$wsHal.Cells.Item(310, 4).Value = 'Alrededor del 77,6% de las personas encuestadas ha recibido formación específica en circo, combinando tradición familiar, espacios comunitarios, autoformación, instituciones técnicas, maestras y maestros particulares y, en algunos casos, procesos internacionales. Existe un consenso amplio sobre la importancia de la formación profesional circense y se priorizan como apoyos clave los recursos económicos, los espacios y condiciones adecuadas de práctica, la formación especializada, el reconocimiento profesional del oficio y el acceso a salud y protección social, delineando una agenda concreta para políticas de profesionalización del sector.'
$wsHal.Cells.Item(311, 1).Value = 109
$wsHal.Cells.Item(311, 2).Value = 1
$wsHal.Cells.Item(311, 3).Value = 'El Festival activa dinámicas económicas y de sosenibilidad cultural'
$wsHal.Cells.Item(311, 4).Value = 'Aunque el apoyo a emprendimientos culturales, la compra de productos o el intercambio con melómanos no constituyen el principal motivo de asistencia en ninguno de los escenarios, sí aparecen de manera consistente en todos ellos como motivaciones complementarias. Este patrón sugiere que el Festival Centro no solo funciona como un espacio de circulación artística, sino también como un dispositivo que contribuye a la sostenibilidad económica del ecosistema cultural, al activar prácticas de consumo cultural, visibilización de agentes y circulación de bienes simbólicos asociados a la música y las artes.
Así, se podría seguir fortaleciendo de manera estratégica los componentes de circulación económica y visibilización de emprendimientos culturales, especialmente en aquellos escenarios con mayor afluencia y diversidad de públicos, sin desdibujar el eje artístico del Festival.
'
$wsHal.Cells.Item(312, 1).Value = 109
$wsHal.Cells.Item(312, 2).Value = 2
$wsHal.Cells.Item(312, 3).Value = 'El Festival no genera mayores impactos negativos en el espacio público'
$wsHal.Cells.Item(312, 4).Value = 'En los cuatro escenarios analizados, la percepción mayoritaria de los asistentes indica que la realización del Festival no modifica sustancialmente problemáticas asociadas al espacio público, como el arrojo de basuras, el parqueo en zonas prohibidas o la contaminación auditiva y visual. 
Por otro lado, en algunos casos, como el Muelle de la FUGA y La Media Torta, se registra una mayor percepción de incremento en la presencia de vendedores informales o en el turismo; sin embargo, estos fenómenos coexisten con una valoración positiva del evento y no se asocian a un deterioro de la convivencia en la zona.
'
$wsHal.Cells.Item(313, 1).Value = 109
$wsHal.Cells.Item(313, 2).Value = 3
$wsHal.Cells.Item(313, 3).Value = 'El Festival articula los desplazamientos, recorridos y consumos culturales en el centro'
$wsHal.Cells.Item(313, 4).Value = 'En todos los escenarios, una proporción mayoritaria de asistentes declaró haber visitado o tener previsto visitar otros espacios del centro antes o después del evento. Los recorridos se concentraron principalmente en equipamientos culturales, cafés, restaurantes, bares, teatros, salas de arte y museos, lo que evidencia una alta capacidad del Festival para articular la oferta cultural, gastronómica y comercial del centro de Bogotá. Este comportamiento refuerza el papel del Festival Centro como dinamizador territorial y como nodo de conexión entre la programación cultural pública y otras actividades económicas y simbólicas del área.'
$wsHal.Cells.Item(314, 1).Value = 109
$wsHal.Cells.Item(314, 2).Value = 4
$wsHal.Cells.Item(314, 3).Value = 'La oferta musical y artística es el principal motivo de asistencia
'
$wsHal.Cells.Item(314, 4).Value = 'En todos los escenarios, los motivos de asistencia se concentran de manera consistente en la presencia de los grupos y artistas, el interés por conocer nuevas propuestas musicales y el reconocimiento previo de algunos de los artistas participantes. El Festival Centro opera simultáneamente como un espacio de encuentro entre artistas y sus audiencias y como una plataforma para el descubrimiento de nuevas propuestas, lo que explica la coexistencia de públicos con trayectoria en el Festival y de personas que asisten por primera vez.'
$wsHal.Cells.Item(315, 1).Value = 110
$wsHal.Cells.Item(315, 2).Value = 1
$wsHal.Cells.Item(315, 3).Value = 'Balance general de la experiencia'
$wsHal.Cells.Item(315, 4).Value = 'En la Noche de Museos y el Festival Patrimonios en Ruana se observa una valoración global muy favorable por parte de los asistentes. Predominan percepciones de alta satisfacción y de calidad en la experiencia cultural, junto con una muy alta disposición a recomendar y a asistir nuevamente en futuras ediciones. En conjunto, los resultados sugieren que ambos formatos funcionan como espacios consistentes para el disfrute del patrimonio y la cultura, con capacidad de sostener la participación en el tiempo.'
$wsHal.Cells.Item(316, 1).Value = 110
$wsHal.Cells.Item(316, 2).Value = 2
$wsHal.Cells.Item(316, 3).Value = 'Públicos, motivaciones y activación cultural'
$wsHal.Cells.Item(316, 4).Value = 'Los hallazgos muestran públicos con trayectorias culturales diversas y motivaciones centradas en hacer algo distinto, explorar la ciudad y acercarse a experiencias culturales específicas. La Noche de Museos destaca por su capacidad de atraer asistentes primerizos, mientras Patrimonios en Ruana combina ese componente con un público recurrente; en ambos casos, la vivencia del evento se asocia con un mayor interés por seguir participando en la oferta cultural de la ciudad.'
$wsHal.Cells.Item(317, 1).Value = 110
$wsHal.Cells.Item(317, 2).Value = 3
$wsHal.Cells.Item(317, 3).Value = 'Encuentro ciudadano, identidad y orgullo por Bogotá D.C.'
$wsHal.Cells.Item(317, 4).Value = 'En ambas mediciones, los eventos aparecen como escenarios de encuentro social con percepciones mayoritarias de respeto e inclusión. Además, se reportan aportes simbólicos claros: conexión con la cultura de la ciudad, fortalecimiento de lazos sociales y comunitarios, y un sentimiento extendido de orgullo por la oferta cultural de Bogotá D.C. En conjunto, los resultados respaldan la idea de que estos eventos no solo convocan público, sino que también fortalecen sentidos compartidos de pertenencia.'
$wsHal.Cells.Item(318, 1).Value = 110
$wsHal.Cells.Item(318, 2).Value = 4
$wsHal.Cells.Item(318, 3).Value = 'Condiciones de acceso, comunicación y oportunidades de mejora'
$wsHal.Cells.Item(318, 4).Value = 'La logística es bien evaluada en aspectos como señalización, apoyo al público y organización general; aun así, el transporte, la movilidad y la seguridad aparecen como factores determinantes para facilitar la asistencia y mejorar la experiencia. En términos de comunicación, Noche de Museos se apoya más en canales digitales, mientras Patrimonios en Ruana depende más del voz a voz, lo que abre una oportunidad para fortalecer la divulgación sin perder el componente comunitario. También se identifican mejoras puntuales en accesibilidad para personas con discapacidad y en temas logísticos asociados al desplazamiento.'
$wsHal.Cells.Item(319, 1).Value = 112
$wsHal.Cells.Item(319, 2).Value = 1
$wsHal.Cells.Item(319, 3).Value = 'Naturaleza de la creación colectiva
'
$wsHal.Cells.Item(319, 4).Value = 'La creación colectiva se caracteriza por su profunda vocación de horizontalidad, solidaridad y libertad. Más que centrarse en la obtención de un producto final, privilegia el proceso como un espacio vivo de exploración, diálogo y construcción conjunta. En este sentido, se configura como un ejercicio humano y político que desafía las jerarquías tradicionales, redistribuye el poder creativo y habilita formas alternativas de relación y producción cultural.'
$wsHal.Cells.Item(320, 1).Value = 112
$wsHal.Cells.Item(320, 2).Value = 2
$wsHal.Cells.Item(320, 3).Value = 'Metodologías, roles y autoría'
$wsHal.Cells.Item(320, 4).Value = 'En los procesos de creación colectiva, las metodologías se construyen desde la participación activa y la co-decisión. Los roles no se imponen de antemano, sino que emergen en función de las necesidades y momentos del proceso, lo que promueve dinámicas más flexibles y colaborativas. El liderazgo se concibe como un acompañamiento empático y orientador, antes que como una instancia de dirección jerárquica. La dramaturgia, por su parte, surge de la investigación, la improvisación y la apertura estructural, y aunque la documentación es una herramienta valiosa, suele ser poco sistemática debido a la naturaleza orgánica del trabajo.'
$wsHal.Cells.Item(321, 1).Value = 112
$wsHal.Cells.Item(321, 2).Value = 3
$wsHal.Cells.Item(321, 3).Value = 'Interdisciplinariedad y articulación de saberes'
$wsHal.Cells.Item(321, 4).Value = 'La creación colectiva se nutre de una integración natural y fluida de múltiples artes, oficios y conocimientos. La apertura metodológica facilita aportes espontáneos entre disciplinas y la incorporación de especialistas externos, enriqueciendo el proceso. La tecnología también se incorpora como recurso creativo y expresivo, privilegiando enfoques de cultura libre y acceso abierto, lo que amplía las posibilidades estéticas y fomenta la innovación en los lenguajes escénicos.'
$wsHal.Cells.Item(322, 1).Value = 112
$wsHal.Cells.Item(322, 2).Value = 4
$wsHal.Cells.Item(322, 3).Value = 'Sostenibilidad e impacto sociocultural'
$wsHal.Cells.Item(322, 4).Value = 'La sostenibilidad de estos procesos depende, ante todo, de los vínculos humanos, la confianza entre los participantes y la claridad de un propósito común. En este modelo, la comunidad no es solo receptora, sino coproductora activa del proceso y sus resultados. Esta dinámica favorece transformaciones significativas tanto en el plano personal como en el territorial, fortaleciendo el sentido de pertenencia, la identidad colectiva y los procesos de memoria y cuidado de los territorios.'
$wsHal.Cells.Item(323, 1).Value = 114
$wsHal.Cells.Item(323, 2).Value = 1
$wsHal.Cells.Item(323, 3).Value = 'La FUGA: Una biblioteca-destino'
$wsHal.Cells.Item(323, 4).Value = 'A diferencia del patrón observado en la mayoría de las bibliotecas públicas —donde predomina el uso para estudio o trabajo—, en la Biblioteca Pública FUGA el principal motivo de visita es recorrerla por primera vez (48,6%). Este valor resulta excepcionalmente alto frente al promedio general del sistema (7,9%) y pone en evidencia su carácter monumental, simbólico y turístico, así como su rol estratégico como equipamiento cultural de alto valor patrimonial en el centro histórico de Bogotá.'
$wsHal.Cells.Item(324, 1).Value = 114
$wsHal.Cells.Item(324, 2).Value = 2
$wsHal.Cells.Item(324, 3).Value = 'Un patrón de afiliación excepcional dentro de la Red'
$wsHal.Cells.Item(324, 4).Value = 'En la Biblioteca Pública Manuel Zapata Olivella – El Tintal, el segundo motivo más frecuente de visita es la afiliación o el uso de otros servicios (18,7%). Este porcentaje es significativamente superior al promedio de la Red, lo que evidencia un desempeño atípico y sugiere que las estrategias de afiliación, fidelización y acercamiento comunitario han tenido un impacto especialmente efectivo en su entorno territorial.'
$wsHal.Cells.Item(325, 1).Value = 114
$wsHal.Cells.Item(325, 2).Value = 3
$wsHal.Cells.Item(325, 3).Value = 'La biblioteca como centro de acceso digital: El caso de Suba'
$wsHal.Cells.Item(325, 4).Value = 'En la Biblioteca Pública Francisco José de Caldas (Suba), el uso de computadores y otros equipos tecnológicos se configura como una motivación de visita altamente relevante, alcanzando el 18,6 % de las respuestas. Esta proporción supera ampliamente el promedio general de las bibliotecas sondeadas (7,0 %), lo que evidencia una demanda territorial específica de acceso a infraestructura tecnológica y posiciona a la biblioteca como un nodo clave para la reducción de brechas digitales en la localidad.'
$wsHal.Cells.Item(326, 1).Value = 114
$wsHal.Cells.Item(326, 2).Value = 4
$wsHal.Cells.Item(326, 3).Value = 'Desconocimiento que frena el uso de los servicios'
$wsHal.Cells.Item(326, 4).Value = 'El principal obstáculo para la vinculación de usuarios potenciales a los servicios bibliotecarios es el desconocimiento de la oferta. Entre las personas que se encuentran en los recintos pero no hacen uso de los servicios, este factor aparece de manera recurrente y con una intensidad superior al promedio del sistema. En la Biblioteca de Usaquén-Servitá, el 44,4 % de estos usuarios señala no conocer la oferta disponible, mientras que en La Peña esta proporción asciende al 50 %, cifras que superan ampliamente el promedio general (20,9 %) y evidencian una brecha crítica de información y comunicación.'
$wsHal.Cells.Item(327, 1).Value = 115
$wsHal.Cells.Item(327, 2).Value = 1
$wsHal.Cells.Item(327, 3).Value = 'Patrones de uso y acompañamiento social'
$wsHal.Cells.Item(327, 4).Value = 'Los resultados muestran que los parques estructurantes priorizados presentan patrones de uso regular, con predominio de visitas semanales y diarias. La asistencia se realiza mayoritariamente en compañía de familiares, y el horario de uso se concentra principalmente en la franja de la mañana. La práctica de deporte o actividad física es la actividad reportada con mayor frecuencia en todos los parques, seguida por usos asociados al acompañamiento y al tránsito cotidiano'
$wsHal.Cells.Item(328, 1).Value = 115
$wsHal.Cells.Item(328, 2).Value = 2
$wsHal.Cells.Item(328, 3).Value = 'Valoración del estado físico y la seguridad'
$wsHal.Cells.Item(328, 4).Value = 'La mayoría de las personas encuestadas califica el estado físico y el mantenimiento de los parques como “bueno” o “excelente”. De igual forma, una alta proporción de usuarios manifiesta sentirse segura al visitar estos espacios. No obstante, una parte de la población reporta percepciones de inseguridad, especialmente asociadas a condiciones como la iluminación, lo que coincide con los aspectos físicos señalados como susceptibles de mejora.
'
$wsHal.Cells.Item(329, 1).Value = 115
$wsHal.Cells.Item(329, 2).Value = 3
$wsHal.Cells.Item(329, 3).Value = 'Reconocimiento comunitario del parque'
$wsHal.Cells.Item(329, 4).Value = 'Los parques son reconocidos por la mayoría de los usuarios como espacios importantes para su comunidad o barrio. Este reconocimiento se presenta de manera consistente en los cinco parques analizados, lo que describe el papel de estos escenarios como referentes del entorno social y del uso cotidiano del espacio público.
'
$wsHal.Cells.Item(330, 1).Value = 115
$wsHal.Cells.Item(330, 2).Value = 4
$wsHal.Cells.Item(330, 3).Value = 'Diferencia entre participación reportada e interés declarado'
$wsHal.Cells.Item(330, 4).Value = 'Los resultados evidencian que la proporción de personas que reporta haber participado en actividades organizadas en los parques es menor que la proporción de usuarios que manifiesta interés en participar en actividades deportivas, recreativas, culturales o comunitarias. Esta diferencia describe una brecha entre la participación efectiva y la disposición declarada por la ciudadanía.
'
$wsHal.Cells.Item(331, 1).Value = 115
$wsHal.Cells.Item(331, 2).Value = 5
$wsHal.Cells.Item(331, 3).Value = 'Barreras y condiciones para el aprovechamiento'
$wsHal.Cells.Item(331, 4).Value = 'La falta de tiempo es la barrera más reportada para el uso y aprovechamiento de los parques, seguida por la percepción de inseguridad y la falta de información sobre la oferta de actividades. Estas barreras se presentan en todos los parques, con variaciones en su magnitud relativa según el escenario.'
$wsHal.Cells.Item(332, 1).Value = 116
$wsHal.Cells.Item(332, 2).Value = 1
$wsHal.Cells.Item(332, 3).Value = 'La motivación existe. 
Las barreras son estructurales y de contexto
'
$wsHal.Cells.Item(332, 4).Value = 'Ciclovía: la mayoría sabe montar bicicleta y reconoce beneficios claros (ahorro de tiempo/dinero, autonomía).

Escuela de la Bici: el principal obstáculo inicial no es el desinterés, sino no saber montar, la inseguridad y la falta de confianza.
'
$wsHal.Cells.Item(333, 1).Value = 116
$wsHal.Cells.Item(333, 2).Value = 2
$wsHal.Cells.Item(333, 3).Value = 'Los programas generan percepción de impacto, pero de forma diferenciada'
$wsHal.Cells.Item(333, 4).Value = 'Ciclovía: alta recurrencia y masividad; impacto percibido en actividad física, motivación semanal y uso recreativo de la bicicleta.

Escuela de la Bici: impacto percibido en aprendizaje, confianza, bienestar emocional y cambio personal, aunque con menor uso posterior por barreras materiales.'
$wsHal.Cells.Item(334, 1).Value = 116
$wsHal.Cells.Item(334, 2).Value = 3
$wsHal.Cells.Item(334, 3).Value = 'La bicicleta transforma la vida cotidiana más que la identidad colectiva'
$wsHal.Cells.Item(334, 4).Value = 'La sociabilidad se expresa principalmente en vínculos cercanos (familia, amigos) más que en colectivos organizados.

El principal beneficio percibido es práctico (ahorro de tiempo/dinero, autonomía), seguido de salud y bienestar.'
$wsHal.Cells.Item(335, 1).Value = 117
$wsHal.Cells.Item(335, 2).Value = 1
$wsHal.Cells.Item(335, 3).Value = 'Comprender la gobernanza deportiva como un sistema complejo, multinivel y en evolución'
$wsHal.Cells.Item(335, 4).Value = 'Los hallazgos de esta investigación muestran que la gobernanza del deporte, la recreación y la actividad física en Bogotá debería entenderse como un sistema complejo en el que convergen actores públicos, privados, comunitarios y académicos con distintos niveles de autoridad, capacidades y responsabilidades. La revisión de literatura confirma que los modelos contemporáneos de gobernanza combinan mecanismos jerárquicos, colaborativos e híbridos, lo cual coincide con la configuración del ecosistema distrital. Esta perspectiva permite reconocer tanto el valor de las estructuras formales del Sistema Nacional del Deporte como la importancia de los vínculos territoriales, las redes sociales y las prácticas comunitarias que sostienen el funcionamiento cotidiano del sistema.'
$wsHal.Cells.Item(336, 1).Value = 117
$wsHal.Cells.Item(336, 2).Value = 2
$wsHal.Cells.Item(336, 3).Value = 'Clarificar el papel de los actores como base para mejorar la articulación y corresponsabilidad'
$wsHal.Cells.Item(336, 4).Value = 'Se considera necesario seguir avanzando en mecanismos que fortalezcan la corresponsabilidad entre niveles de gobierno, organizaciones deportivas, sector privado, ciudadanía y academia. La claridad en los roles y la comprensión de las relaciones entre actores son condiciones esenciales para que la política deportiva distrital evolucione hacia formas más consistentes de cooperación, planificación conjunta y toma de decisiones basada en evidencia, la tipología de actores es un punto de partida para esta definición.'
$wsHal.Cells.Item(337, 1).Value = 117
$wsHal.Cells.Item(337, 2).Value = 3
$wsHal.Cells.Item(337, 3).Value = 'Aplicar la medición para caracterizar a los actores y fortalecer el sistema de información'
$wsHal.Cells.Item(337, 4).Value = 'El desarrollo de las categorías analíticas, del formulario cuantitativo y de la propuesta metodológica cualitativa constituye un insumo clave para avanzar hacia un sistema de información más sólido y útil para la toma de decisiones. La aplicación de estos instrumentos permitirá caracterizar de manera rigurosa a los actores del ecosistema deportivo y facilitará una comprensión más completa, contextualizada y accionable sobre el funcionamiento real del ecosistema deportivo y sobre las oportunidades para fortalecer su gobernanza.'
$wsHal.Cells.Item(338, 1).Value = 117
$wsHal.Cells.Item(338, 2).Value = 4
$wsHal.Cells.Item(338, 3).Value = 'Usar el instrumento como herramienta progresiva para la mejora de la política pública'
$wsHal.Cells.Item(338, 4).Value = 'Se recomienda concebir el instrumento desarrollado no como un ejercicio aislado de levantamiento de información, sino como una herramienta estratégica de uso progresivo para el fortalecimiento de la política pública del deporte en Bogotá. Su aplicación periódica y su articulación con otros sistemas de información permitirán identificar brechas, monitorear avances en la articulación institucional y ajustar intervenciones de manera gradual. En este sentido, el instrumento ofrece una base para avanzar hacia procesos de planeación, seguimiento y evaluación más coherentes, orientados al aprendizaje institucional y a la generación de valor público en el ecosistema deportivo.

En este marco, el instrumento constituye una base para promover procesos de aprendizaje institucional, mejorar la coordinación entre actores y fortalecer la generación de valor público en el ecosistema deportivo, a partir de decisiones informadas y basadas en evidencia.
'
$wsHal.Cells.Item(339, 1).Value = 119
$wsHal.Cells.Item(339, 2).Value = 1
$wsHal.Cells.Item(339, 3).Value = 'Balance general y experiencia del público'
$wsHal.Cells.Item(339, 4).Value = 'El Festival de Verano 2025 presenta una valoración mayoritariamente positiva en Bogotá D.C., con satisfacción frente a la programación y a la experiencia general de participación. Como recomendación, se sugiere mantener temas logísticos que más inciden en el disfrute de la experiencia como personal de apoyo, información clara, entre otros, especialmente en escenarios de alta afluencia. '
$wsHal.Cells.Item(340, 1).Value = 119
$wsHal.Cells.Item(340, 2).Value = 2
$wsHal.Cells.Item(340, 3).Value = 'Accesibilidad, movilidad y seguridad como factores críticos'
$wsHal.Cells.Item(340, 4).Value = 'La evidencia indica que la experiencia no depende solo de las actividades, sino de condiciones de acceso, circulación y permanencia. Se recomienda fortalecer la gestión de movilidad (incluida la salida) y de seguridad, con medidas preventivas y coordinación interinstitucional, priorizando los puntos y franjas horarias donde se concentran los mayores flujos. '
$wsHal.Cells.Item(341, 1).Value = 119
$wsHal.Cells.Item(341, 2).Value = 3
$wsHal.Cells.Item(341, 3).Value = 'Convivencia y gestión del entorno territorial'
$wsHal.Cells.Item(341, 4).Value = 'En términos de convivencia, el festival es percibido como espacio favorable para el encuentro, pero la medición a vecinos y ciudadanía muestra tensiones asociadas al entorno (parqueo indebido, residuos, ruido y otras dinámicas percibidas). Se recomienda reforzar estrategias de ordenamiento, mitigación ambiental y corresponsabilidad ciudadana en zonas de influencia, para equilibrar disfrute y bienestar barrial. '
$wsHal.Cells.Item(342, 1).Value = 119
$wsHal.Cells.Item(342, 2).Value = 4
$wsHal.Cells.Item(342, 3).Value = 'Dinámica económica local y formalización de oportunidades'
$wsHal.Cells.Item(342, 4).Value = 'Los resultados sugieren una activación económica asociada al festival (ingresos reportados por vendedores informales y participación de emprendimientos), aunque con comportamientos heterogéneos. Se recomienda consolidar lineamientos de participación, articulación con emprendimientos y vendedores, y mecanismos de seguimiento que permitan aprovechar mejor el potencial económico sin perder el enfoque cultural, recreativo y de espacio público del evento.'
$wsHal.Cells.Item(343, 1).Value = 120
$wsHal.Cells.Item(343, 2).Value = 1
$wsHal.Cells.Item(343, 3).Value = 'Balance general de satisfacción'
$wsHal.Cells.Item(343, 4).Value = 'La investigación consolida resultados de la Bienal Internacional de Arte y Ciudad BOG25, el Concurso Internacional de Violín 2025 y Navidad es Cultura 2025, donde se identifica una experiencia ampliamente positiva, con altos niveles de recomendación, orgullo y disposición a regresar en futuras ediciones, lo que sugiere una aceptación sólida de los formatos en Bogotá D.C. '
$wsHal.Cells.Item(344, 1).Value = 120
$wsHal.Cells.Item(344, 2).Value = 2
$wsHal.Cells.Item(344, 3).Value = 'Aportes culturales y formativos para los públicos'
$wsHal.Cells.Item(344, 4).Value = 'De forma agregada, los resultados muestran que la experiencia no se limita al disfrute: una proporción alta reporta aprendizajes, descubrimiento y reflexión, así como un incremento en el interés por seguir participando en la oferta cultural. En este sentido, se destaca el descubrimiento de nuevos artistas o discursos, la reflexión provocada por las obras y la motivación para asistir a otros eventos; además, se observa un fortalecimiento del interés por la música clásica y la formación de públicos en Bogotá D.C. '
$wsHal.Cells.Item(345, 1).Value = 120
$wsHal.Cells.Item(345, 2).Value = 3
$wsHal.Cells.Item(345, 3).Value = 'Encuentro ciudadano, identidad y orgullo por Bogotá D.C.'
$wsHal.Cells.Item(345, 4).Value = 'En conjunto, los eventos se perciben como espacios que favorecen el encuentro social, la convivencia y el fortalecimiento de sentidos de identidad y pertenencia. Los hallazgos también muestran un componente simbólico fuerte de orgullo por Bogotá D.C. y una lectura del arte y la cultura como herramientas que pueden reunir a personas diversas y abrir conversaciones. Adicionalmente, tras la experiencia, una parte importante de asistentes reporta mayor confianza en la articulación institucional y en el uso transparente de recursos públicos. '
$wsHal.Cells.Item(346, 1).Value = 120
$wsHal.Cells.Item(346, 2).Value = 4
$wsHal.Cells.Item(346, 3).Value = 'Acceso y comunicación: fortalezas y oportunidades de mejora'
$wsHal.Cells.Item(346, 4).Value = 'En términos generales, se reportan condiciones de participación mayoritariamente favorables, con énfasis en la accesibilidad física y la logística; sin embargo, los resultados sugieren un reto claro en comunicación previa (cuando una proporción de asistentes percibe información poco clara o confusa) y, en eventos desarrollados en espacio público, la mediación puede reforzarse para que más personas identifiquen que lo que observan hace parte del evento. En síntesis, el balance operativo es positivo, pero hay margen para fortalecer divulgación y mediación para ampliar alcance y mejorar la experiencia.'
$wsHal.Cells.Item(347, 1).Value = 121
$wsHal.Cells.Item(347, 2).Value = 4
$wsHal.Cells.Item(347, 3).Value = 'Mucho orgullo, confianza aún frágil'
$wsHal.Cells.Item(347, 4).Value = 'La encuesta muestra un considerable orgullo por Bogotá, pero una confianza interpersonal limitada: la mayoría dice poder confiar en menos de la mitad de las personas. Esta brecha sugiere que el sentido de pertenencia no se traduce automáticamente en cooperación cotidiana. Para cerrarla, los laboratorios deben priorizar mediaciones comunitarias, pedagogía de cuidado del espacio público y acciones que incrementen encuentros repetidos y seguros entre vecinos. Deporte, juego y arte son palancas idóneas para “hacer juntos” con bajas barreras de entrada. Medir confianza barrial antes y después de cada ciclo permitirá verificar si las intervenciones logran convertir orgullo en vínculos efectivos.'
$wsHal.Cells.Item(348, 1).Value = 121
$wsHal.Cells.Item(348, 2).Value = 5
$wsHal.Cells.Item(348, 3).Value = 'Gobernanza compartida con palancas culturales'
$wsHal.Cells.Item(348, 4).Value = 'Los hallazgos apuntan a corresponsabilidad distribuida: fuerza pública, entidades distritales, vecinos y JAL son percibidos como actores clave según el tipo de problema. Esta lectura respalda una gobernanza compartida, con roles claros y coordinación por barrio, donde cultura ciudadana, bibliotecas y arte urbano actúan como catalizadores. La priorización debe enfocarse en inseguridad y basuras, conectando intervenciones de seguridad situacional con cuidado del entorno y gestión de residuos. Diseñar participación flexible (horarios extendidos, micro-tareas, reconocimiento simbólico) reduce barreras para ese tercio que no puede involucrarse. Con tableros de seguimiento por barrio, la estrategia podrá iterar con evidencia y escalar lo que funciona.'
$wsHal.Cells.Item(349, 1).Value = 128
$wsHal.Cells.Item(349, 2).Value = 1
$wsHal.Cells.Item(349, 3).Value = 'Un instrumento que lee el oficio como economía cultural'
$wsHal.Cells.Item(349, 4).Value = 'El nuevo formulario no solo recoge información: reconstruye la lógica productiva del oficio artesanal. Aunque mantiene seis ejes estructurales (sociodemografía, identificación, elaboración, comercialización, valor patrimonial y consentimiento), incorpora variables que permiten entender al artesano como unidad productiva inserta en dinámicas territoriales y de mercado.'
$wsHal.Cells.Item(350, 1).Value = 128
$wsHal.Cells.Item(350, 2).Value = 2
$wsHal.Cells.Item(350, 3).Value = 'Georreferenciación: el espacio público como territorio productivo'
$wsHal.Cells.Item(350, 4).Value = 'La incorporación de georreferenciación de puntos de elaboración y venta, junto con la identificación de tipos de espacios (públicos y privados), convierte el instrumento en una herramienta de análisis territorial. Esta información permite mapear las dinámicas reales de ocupación económica del espacio público, identificar concentraciones, circuitos y patrones de uso, y generar insumos estratégicos para la gestión urbana y la toma de decisiones sobre regulación y permisos.'
$wsHal.Cells.Item(351, 1).Value = 128
$wsHal.Cells.Item(351, 2).Value = 3
$wsHal.Cells.Item(351, 3).Value = 'Hacedores, creadores y comercializadores: una tipología para entender el rol de los artesanos en la cadena de valor'
$wsHal.Cells.Item(351, 4).Value = 'La introducción de preguntas sobre rol principal frente a las piezas, etapas del proceso en las que se participa y vínculo con la venta permitió construir una clasificación en cuatro categorías: hacedor/a, hacedor/a‑comercializador/a, creador/a y comercializador/a. Esta tipología, basada en una definición de artesano que articula medio de vida, expresión cultural y dominio técnico, ofrece un insumo clave para orientar estrategias diferenciadas de formación, visibilización y fortalecimiento según el lugar que cada persona ocupa en la cadena de valor artesanal'
$wsHal.Cells.Item(352, 1).Value = 133
$wsHal.Cells.Item(352, 2).Value = 1
$wsHal.Cells.Item(352, 3).Value = 'La violencia sexual no es un desvío, es un mandato aprendido'
$wsHal.Cells.Item(352, 4).Value = 'La investigación evidencia que la violencia sexual no aparece como un hecho aislado o accidental, sino como una expresión coherente de la masculinidad hegemónica aprendida en distintos entornos sociales. El ejercicio de la violencia se asocia a expectativas de poder, control, dominación y validación identitaria masculina. Estas prácticas son reforzadas por la familia, los pares, los medios y, en ocasiones, por la ausencia de una educación sexual integral. La violencia se convierte así en un mecanismo para “demostrar hombría”, más que en una desviación individual.'
$wsHal.Cells.Item(353, 1).Value = 133
$wsHal.Cells.Item(353, 2).Value = 2
$wsHal.Cells.Item(353, 3).Value = 'El hogar es el principal escenario de riesgo y silenciamiento'
$wsHal.Cells.Item(353, 4).Value = 'Una proporción significativa de los hechos de violencia sexual ocurre en la vivienda y es ejercida por personas cercanas a la víctima. El entorno familiar, lejos de ser siempre un espacio protector, puede funcionar como un escenario de encubrimiento, negación o minimización del daño. Esta dinámica dificulta la denuncia, la intervención temprana y la activación de rutas de protección. La naturalización del abuso dentro del hogar refuerza el silencio y perpetúa la violencia.'
$wsHal.Cells.Item(354, 1).Value = 133
$wsHal.Cells.Item(354, 2).Value = 3
$wsHal.Cells.Item(354, 3).Value = 'Negar, minimizar y justificar: la principal barrera para la responsabilidad'
$wsHal.Cells.Item(354, 4).Value = 'Tanto los ofensores como sus familias tienden a negar, minimizar o justificar la violencia sexual, lo que constituye una barrera estructural para los procesos de responsabilización y reparación. Estas narrativas desplazan la culpa hacia las víctimas, apelan a la cercanía afectiva o normalizan el abuso como “error” o “juego”. Esta falta de reconocimiento del daño limita la efectividad de las sanciones y de los procesos restaurativos. Sin responsabilidad no hay transformación posible.'
$wsHal.Cells.Item(355, 1).Value = 133
$wsHal.Cells.Item(355, 2).Value = 4
$wsHal.Cells.Item(355, 3).Value = 'La sexualidad se aprende sin consentimiento ni empatía'
$wsHal.Cells.Item(355, 4).Value = 'Los hallazgos muestran que muchos jóvenes adquieren conocimientos sobre sexualidad a través de pares, pornografía y medios, sin mediación crítica ni enfoque ético. Esto produce ideas distorsionadas sobre el consentimiento, la reciprocidad y el deseo, normalizando la cosificación y la dominación. La ausencia de educación sexual integral favorece prácticas sexuales violentas y dificulta el desarrollo de empatía hacia las víctimas. La sexualidad se configura más como ejercicio de poder que como relación.'
$wsHal.Cells.Item(356, 1).Value = 133
$wsHal.Cells.Item(356, 2).Value = 5
$wsHal.Cells.Item(356, 3).Value = 'Transformar las masculinidades es clave para prevenir la violencia'
$wsHal.Cells.Item(356, 4).Value = 'La erradicación de la violencia sexual requiere intervenir directamente en la construcción social de las masculinidades. El estudio muestra que trabajar únicamente sobre las víctimas o desde una lógica punitiva es insuficiente. Se necesitan estrategias preventivas, educativas, culturales y restaurativas que involucren activamente a los hombres, desde la infancia hasta la adultez. Promover masculinidades no violentas es una condición central para la justicia, la convivencia y la transformación cultural.'

Write-Output "done"
